$wb = $excel.ActiveWorkbook

# Re-introduce the base fee shift that was lost along the way during updates.
$wsInternal = $wb.Worksheets.Item("Outputs_Internal")
$wsInternal.Range("H7").Value = 0.075

# Update the recorded selection on Outputs_Internal.
$wsInternal.Range("C11").Select()

# Update the recorded selection/scroll position on Prices.
$wsPrices = $wb.Worksheets.Item("Prices")
$wsPrices.Activate()
$wsPrices.Range("D24").Select()

# Restore Inputs as the active sheet (unchanged tab selection in the edit).
$wsInputs = $wb.Worksheets.Item("Inputs")
$wsInputs.Activate()
